$d = $word.ActiveDocument

# Remove the (unused) magenta highlight from the two respiratory-insufficiency
# symptom entries — "Бронхо спазам" and "Пнеумонија" — covering both the run
# text and the paragraph mark itself so no <w:highlight> survives on either.
$targets = @("Бронхо спазам", "Пнеумонија")

foreach ($p in $d.Paragraphs) {
    foreach ($t in $targets) {
        if ($p.Range.Text -like "*$t*") {
            $p.Range.Font.HighlightColorIndex = 0
        }
    }
}

Write-Output "done"
